$wb = $excel.ActiveWorkbook

# --- "Settings" sheet (sheet1): insert two new config rows right after the header row ---
$ws1 = $wb.Worksheets.Item("Settings")

# Insert 2 blank rows above the current row 2 (shifts everything else down by 2)
$ws1.Rows("2:3").Insert()

# New Name/Value pairs for the Orchestrator queue settings used by the
# ConocoPhillips Dispatcher bot
$ws1.Range("A2").Value = "OrchestratorQueueName"
$ws1.Range("A3").Value = "OrchestratorQueueFolder"
$ws1.Range("B2").Value = "ConocoPhillips_Files"
$ws1.Range("B3").Value = "Shared"

# The inserted rows copy the bold header formatting from row 1 by default;
# clear that back down to a plain/default look for the two label/value cells
# (matches the rest of the Name/Value rows in the sheet).
$ws1.Range("A2:B3").ClearFormats()

# --- "Constants" sheet (sheet2): cursor moved, no content changes ---
$ws2 = $wb.Worksheets.Item("Constants")
$ws2.Range("D9").Select()

# Leave "Settings" as the active sheet/selection, matching the saved cursor
# position after the edits were made.
$ws1.Range("C21").Select()
